# daily auto push: 2026-02-23 19:19 UTC
# Insert two new daily-log rows (2026/02/23 and 2026/02/24) above the existing
# "2026/12/29" block, shifting all rows from the old 869 onward down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 869 (old row 869 and everything after it
# shifts down to 871.. / 912).
$ws.Rows("869:870").Insert()

# Force column A on the two new rows to be stored as plain text so the
# "yyyy/mm/dd"-looking strings are not auto-converted to date serial values
# (matching the existing inlineStr/text date cells used throughout the sheet).
$ws.Range("A869:A870").NumberFormat = "@"

# New row 869: 2026/02/23 (Mon)
$ws.Cells.Item(869, 1).Value = "2026/02/23"
$ws.Cells.Item(869, 2).Value = "月"
$ws.Cells.Item(869, 3).Value = 23
$ws.Cells.Item(869, 4).Value = 201

# New row 870: 2026/02/24 (Tue)
$ws.Cells.Item(870, 1).Value = "2026/02/24"
$ws.Cells.Item(870, 2).Value = "火"
$ws.Cells.Item(870, 3).Value = 2
$ws.Cells.Item(870, 4).Value = 201
